$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting rows 71:139 down to 72:140.
$ws.Range("A71").EntireRow.Insert()

# Populate the newly inserted row 71 with the new data point.
$ws.Cells.Item(71, 1).Value = 10
$ws.Cells.Item(71, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(71, 3).Value = "La Araucanía"
$ws.Cells.Item(71, 4).Value = 44484
$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
$ws.Cells.Item(71, 5).Value = 9
$ws.Cells.Item(71, 6).Value = 100112005
$ws.Cells.Item(71, 7).Value = "Puerro"
$ws.Cells.Item(71, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 50
$ws.Cells.Item(71, 11).Value = 6500
$ws.Cells.Item(71, 12).Value = 7000
$ws.Cells.Item(71, 13).Value = 6800
$ws.Cells.Item(71, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(71, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(71, 16).Value = 567
$ws.Cells.Item(71, 17).Value = 12
$ws.Cells.Item(71, 18).Value = "Hortaliza"
